$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting A:E to B:F
$ws.Columns.Item(1).Insert()

# New header for the ID column
$ws.Cells.Item(1, 1).Value = "ID"

# Match the header formatting (bold, centered, bordered) used by the other
# header cells (B1:F1) so A1 looks consistent with the rest of row 1.
$ws.Cells.Item(1, 2).Copy() | Out-Null
$ws.Cells.Item(1, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# ID values for each data row (rows 2-25)
$ids = @(
    "Hb 2",
    "Hb 3",
    "S 24",
    "S 28",
    "Hb 107",
    "Hb 66",
    "Hb 69",
    "Hb 95",
    "Hb 99",
    "Hb 92",
    "Hb 40",
    "Hb 41",
    "S 11",
    "Hb 57",
    "S 21",
    "S 22",
    "S 3",
    "S 4",
    "S 5",
    "Hb 74",
    "Hb 79",
    "Hb 32",
    "S 15",
    "S 16"
)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}

Write-Output "Done"
